$wb = $excel.ActiveWorkbook

# Header row is identical (lower-cased) across all 15 sheets
$header = @{
    'A1' = 'ratings'
    'B1' = 'api_and_integration_support'
    'C1' = 'pricing_details'
    'D1' = 'deployment_support'
    'E1' = 'customer_support_options'
    'F1' = 'training_platforms'
    'G1' = 'vendor_details'
    'H1' = 'features'
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($col in $header.Keys) {
        $ws.Range($col).Value = $header[$col]
    }
}

# Sheet index 1: {"company_name":"ServiceNow","year_founded":2004,"country":"United States"}
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = '{"total_reviews":1000,"ease_of_use":4.5,"features":4.7,"design":4.3,"support":4.2,"overall":4.6}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":true,"BigID":true,"Cozyroc SSIS+ Suite":true,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":true,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":true,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"ServiceNow","year_founded":2004,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","IT Operations Management","Cloud Management","Security Management"]'

# Sheet index 2: {"company_name":"SolarWinds","year_founded":1999,"country":"United States"}
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = '{"total_reviews":500,"ease_of_use":4.2,"features":4,"design":3.8,"support":4,"overall":4.1}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":true,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"SolarWinds","year_founded":1999,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 3: {"company_name":"ManageEngine","year_founded":1999,"country":"India"}
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = '{"total_reviews":300,"ease_of_use":4,"features":3.8,"design":3.5,"support":3.7,"overall":3.9}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":true,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"ManageEngine","year_founded":1999,"country":"India"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 4: {"company_name":"TOPdesk","year_founded":1993,"country":"Netherlands"}
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = '{"total_reviews":200,"ease_of_use":4.1,"features":3.9,"design":3.7,"support":3.8,"overall":3.9}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"TOPdesk","year_founded":1993,"country":"Netherlands"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 5: {"company_name":"SymphonyAI","year_founded":2017,"country":"United States"}
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = '{"total_reviews":150,"ease_of_use":4.3,"features":4.5,"design":4,"support":3.9,"overall":4.2}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":true,"BigID":true,"Cozyroc SSIS+ Suite":true,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":true,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":true,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"SymphonyAI","year_founded":2017,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","IT Operations Management","Cloud Management","Security Management","AI-powered automation"]'

# Sheet index 6: {"company_name":"Atlassian","year_founded":2002,"country":"Australia"}
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = '{"total_reviews":800,"ease_of_use":4.4,"features":4.3,"design":4,"support":4.1,"overall":4.3}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":true,"BigID":false,"Cozyroc SSIS+ Suite":true,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":true,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"Atlassian","year_founded":2002,"country":"Australia"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","Project Management","Agile Development"]'

# Sheet index 7: {"company_name":"Cherwell Software","year_founded":2000,"country":"United States"}
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = '{"total_reviews":100,"ease_of_use":3.8,"features":3.7,"design":3.5,"support":3.6,"overall":3.7}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"Cherwell Software","year_founded":2000,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 8: {"company_name":"Freshworks","year_founded":2010,"country":"India"}
$ws = $wb.Worksheets.Item(8)
$ws.Range("A2").Value = '{"total_reviews":400,"ease_of_use":4.2,"features":4,"design":3.9,"support":4,"overall":4.1}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":true,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"Freshworks","year_founded":2010,"country":"India"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","IT Operations Management"]'

# Sheet index 9: {"company_name":"SysAid Technologies","year_founded":1999,"country":"Israel"}
$ws = $wb.Worksheets.Item(9)
$ws.Range("A2").Value = '{"total_reviews":250,"ease_of_use":3.9,"features":3.7,"design":3.6,"support":3.8,"overall":3.8}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":true,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"SysAid Technologies","year_founded":1999,"country":"Israel"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","Remote Support"]'

# Sheet index 10: {"company_name":"BMC Software","year_founded":1980,"country":"United States"}
$ws = $wb.Worksheets.Item(10)
$ws.Range("A2").Value = '{"total_reviews":75,"ease_of_use":3.5,"features":3.6,"design":3.3,"support":3.4,"overall":3.5}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"BMC Software","year_founded":1980,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 11: {"company_name":"Ivanti","year_founded":1994,"country":"United States"}
$ws = $wb.Worksheets.Item(11)
$ws.Range("A2").Value = '{"total_reviews":125,"ease_of_use":4,"features":3.8,"design":3.7,"support":3.6,"overall":3.8}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"Ivanti","year_founded":1994,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","IT Operations Management"]'

# Sheet index 12: {"company_name":"EV Technologies","year_founded":1998,"country":"United Kingdom"}
$ws = $wb.Worksheets.Item(12)
$ws.Range("A2").Value = '{"total_reviews":100,"ease_of_use":3.7,"features":3.6,"design":3.4,"support":3.5,"overall":3.6}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"EV Technologies","year_founded":1998,"country":"United Kingdom"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 13: {"company_name":"SolarWinds","year_founded":1999,"country":"United States"}
$ws = $wb.Worksheets.Item(13)
$ws.Range("A2").Value = '{"total_reviews":200,"ease_of_use":3.8,"features":3.6,"design":3.5,"support":3.7,"overall":3.7}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"SolarWinds","year_founded":1999,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics"]'

# Sheet index 14: {"company_name":"TeamDynamix","year_founded":2001,"country":"United States"}
$ws = $wb.Worksheets.Item(14)
$ws.Range("A2").Value = '{"total_reviews":150,"ease_of_use":3.9,"features":3.7,"design":3.6,"support":3.8,"overall":3.8}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":false,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":true}'
$ws.Range("G2").Value = '{"company_name":"TeamDynamix","year_founded":2001,"country":"United States"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","Project Management"]'

# Sheet index 15: {"company_name":"InvGate","year_founded":2005,"country":"Argentina"}
$ws = $wb.Worksheets.Item(15)
$ws.Range("A2").Value = '{"total_reviews":100,"ease_of_use":4,"features":3.8,"design":3.7,"support":3.9,"overall":3.9}'
$ws.Range("B2").Value = '{"api_access":true,"integration_support":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other available integrations":true}}'
$ws.Range("C2").Value = '{"pricing_tiers":"Contact vendor for pricing","free_version_availability":true,"free_trial_availability":true}'
$ws.Range("D2").Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range("E2").Value = '{"phone_support":true,"24/7_live_support":false,"online_support":true}'
$ws.Range("F2").Value = '{"documentation":true,"webinars":true,"live_online_sessions":false,"in_person_training":false}'
$ws.Range("G2").Value = '{"company_name":"InvGate","year_founded":2005,"country":"Argentina"}'
$ws.Range("H2").Value = '["Incident Management","Problem Management","Change Management","Asset Management","Knowledge Management","Service Catalog","Request Management","Service Level Management","Reporting and Analytics","IT Operations Management"]'
